# Insert a new data row at row 26 (weekly Berenjena price entry), pushing
# the existing rows 26:68 down to 27:69, then populate the new row with
# the reported values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 26..68 down to 27..69
$ws.Rows("26:26").Insert()

# Populate the newly inserted row 26 with the new observation
$ws.Range("A26").Value = 5
$ws.Range("B26").Value = "Macroferia Regional de Talca"
$ws.Range("C26").Value = "Maule"
$ws.Range("D26").Value = 44469
$ws.Range("E26").Value = 7
$ws.Range("F26").Value = 100112001
$ws.Range("G26").Value = "Berenjena"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Segunda"
$ws.Range("J26").Value = 300
$ws.Range("K26").Value = 8000
$ws.Range("L26").Value = 8000
$ws.Range("M26").Value = 8000
$ws.Range("N26").Value = "`$/caja 50 unidades"
$ws.Range("O26").Value = "Región de Arica y Parinacota"
$ws.Range("P26").Value = 160
$ws.Range("Q26").Value = 50
$ws.Range("R26").Value = "Hortaliza"
